# Add a new "time_taken" metadata column (F) to the panel worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1: copy formatting from the neighbouring header cell (E1)
# so it picks up the same bold/border/center style, then set its text.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Per-row "time_taken" values (column F, rows 2-11), matching the plain
# (unstyled) formatting already used by the other data columns.
$times = @(
    "2021-10-05 10:52:41.190879",
    "2021-10-05 10:52:41.190890",
    "2021-10-05 10:52:41.190894",
    "2021-10-05 10:52:41.190897",
    "2021-10-05 10:52:41.190901",
    "2021-10-05 10:52:41.190904",
    "2021-10-05 10:52:41.190907",
    "2021-10-05 10:52:41.190910",
    "2021-10-05 10:52:41.190913",
    "2021-10-05 10:52:41.190916"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
